# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# 236312fd-... and 97cee9a5-... files have been handed back (targets +
# handback files generated), and refreshes the "Status" column text and
# the relevant column widths so the new file-name columns are readable.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/23ff3d8781311e723d4c93c1661f540a794fa60d/e2e/"

$mdName1 = "236312fd-583d-4ff2-a98f-82bb1ee12462.md"
$mdName2 = "97cee9a5-f440-48b0-9ec8-2d136489f0ec.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de columns (E, F) and refresh the
# status text shown there (it mirrors the per-language Status column).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): record the handback, i.e. populate
# "Latest Target File" (I) and "Latest Handback File" (J), refresh
# "Latest Handback DateTime" (K) and the Status column (C), widen the
# columns that now hold long file names, and add hyperlinks on the new
# target-file cells (mirroring the existing source-file hyperlinks).
# ---------------------------------------------------------------------

$langs = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-09-06 00:28:49" },
    @{ Name = "de-de"; HandbackTime = "2016-09-06 00:28:57" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Widen Status (C) and the two file-name columns (I, J)
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40

    # Refresh Status column
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Populate Latest Target File (I) and Latest Handback File (J). The
    # handback file matches the already-generated handoff xliff (G) for
    # the same row.
    $ws.Range("I2").Value = $mdName1
    $ws.Range("J2").Value = $ws.Range("G2").Value()

    $ws.Range("I3").Value = $mdName2
    $ws.Range("J3").Value = $ws.Range("G3").Value()

    # Refresh Latest Handback DateTime (K)
    $ws.Range("K2").Value = $lang.HandbackTime
    $ws.Range("K3").Value = $lang.HandbackTime

    # Rebuild hyperlinks so the new "Latest Target File" cells link out to
    # the same source files as column A, interleaved in row order.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + $mdName1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $baseUrl + $mdName1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + $mdName2, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $baseUrl + $mdName2, "", "", $mdName2)
}
